# Experiment_results.xlsx — "Answered all questions and implemented all code"
#
# - Re-ran the last experiment (SGDClassifier / LogisticRegression block,
#   rows 38-40) with updated scores.
# - Cleared the yellow highlight on the AdaBoostClassifier block
#   (rows 10-12, cols A-F) back to white.
# - Moved the active selection to C18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the yellow fill on the AdaBoostClassifier rows (10-12) -> white ---
# ThemeColor 2 == xlThemeColorLight1 (Background 1 / white), which is what
# Excel itself writes out as <fgColor theme="0"/> in styles.xml.
$ws.Range("A10:F12").Interior.ThemeColor = 2

# --- Update the re-run results for the last classifier (rows 38-40) ---
$ws.Range("E38").Value = 0.8593
$ws.Range("F38").Value = 0.7647

$ws.Range("C39").Value = 0.0034
$ws.Range("D39").Value = 0.0002
$ws.Range("E39").Value = 0.8562
$ws.Range("F39").Value = 0.7914

$ws.Range("C40").Value = 0.0028
$ws.Range("D40").Value = 0.0001
$ws.Range("E40").Value = 0.8468
$ws.Range("F40").Value = 0.806

# --- Move the selection, matching the author's final cursor position ---
$ws.Range("C18").Select()
